$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells with same style as the other header cells (copy style from AC1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerSrc = $ws.Range("AC1")
$headerDst = $ws.Range("AD1:AF1")
$headerSrc.Copy()
$headerDst.PasteSpecial(-4122)

# Fill in team record data for each data row (rows 2 through 58)
for ($r = 2; $r -le 58; $r++) {
    $ws.Cells.Item($r, 30).Value = 75
    $ws.Cells.Item($r, 31).Value = 87
    $ws.Cells.Item($r, 32).Value = 0
}
